$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Truncate the coordinate values in Q2 (Ost) and R2 (Nord) to integers
$ws.Range("Q2").Value = 575785
$ws.Range("R2").Value = 6300743

# Remove the Starttid (Z2) and Sluttid (AB2) values entirely
$ws.Range("Z2").ClearContents()
$ws.Range("AB2").ClearContents()
